$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data (GitHub Actions scheduled refresh).
# Column D ("Price") values that look numeric are assigned with a leading
# apostrophe so Excel keeps them as literal text (preserving trailing
# zeros/formatting exactly as scraped) instead of auto-coercing to a number.
$ws.Range("D2").Value = "68.354.96"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.774.78"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'596.04"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'168.40"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "3.770.77"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").Value = "'36.44"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "4.409.59"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "3.778.12"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "68.332.44"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'18.20"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("D19").Value = "'7.05"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'11.00"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").Value = "'469.01"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'0.702"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "'84.80"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "'12.19"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'10.22"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "3.922.72"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "'30.07"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'9.27"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D37").Value = "3.729.08"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E39").Value = "  -10.21%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.308"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").Value = "'43.93"
$ws.Range("E46").Value = "  +12.42%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'8.58"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "'406.76"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "'45.57"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'145.78"
$ws.Range("E51").Value = "  +2.52%  "

# Reset number formatting on the text-forced Price cells back to the
# workbook's default (General/no style) so only the value changes - the
# apostrophe-prefix entry above would otherwise mark the cell "quote
# prefixed" and give it a distinct style.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
